$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark cells whose new numeric-looking value must stay text
# (Excel auto-converts clean numeric strings to numbers on assignment;
#  pre-formatting as Text keeps these as strings, matching the source data).
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D11",
    "D15",
    "D19",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D45",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Step 2: write the updated cell values
$ws.Range('D2').Value = '61.284.92'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '2.680.63'
$ws.Range('E3').Value = '  +2.38%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '534.43'
$ws.Range('E5').Value = '  +3.96%  '
$ws.Range('D6').Value = '157.31'
$ws.Range('E6').Value = '  +1.88%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('D9').Value = '6.61'
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('E10').Value = '  +4.77%  '
$ws.Range('D11').Value = '0.356'
$ws.Range('E11').Value = '  +2.67%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '3.147.35'
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = '61.337.57'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = '22.26'
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = '2.677.17'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('D19').Value = '357.22'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '10.79'
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('E21').Value = '  +3.43%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '61.97'
$ws.Range('E23').Value = '  +1.77%  '
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('D25').Value = '0.169'
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('D27').Value = '0.0₃0872'
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('D28').Value = '7.45'
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '6.20'
$ws.Range('E30').Value = '  +5.65%  '
$ws.Range('D31').Value = '19.66'
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('E32').Value = '  +3.51%  '
$ws.Range('D33').Value = '150.34'
$ws.Range('E33').Value = '  -1.15%  '
$ws.Range('D34').Value = '4.18'
$ws.Range('E34').Value = '  +4.90%  '
$ws.Range('E35').Value = '  +2.32%  '
$ws.Range('D36').Value = '0.923'
$ws.Range('E36').Value = '  +9.33%  '
$ws.Range('D37').Value = '0.892'
$ws.Range('E37').Value = '  +2.65%  '
$ws.Range('D38').Value = '1.51'
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '308.84'
$ws.Range('E39').Value = '  +5.17%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '36.98'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').Value = '0.655'
$ws.Range('E42').Value = '  +4.76%  '
$ws.Range('D43').Value = '20.99'
$ws.Range('E43').Value = '  +5.98%  '
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').Value = '0.0570'
$ws.Range('E45').Value = '  +2.75%  '
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('E47').Value = '  +2.68%  '
$ws.Range('D48').Value = '0.0242'
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('D49').Value = '19.28'
$ws.Range('E49').Value = '  +8.48%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '10.35'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.009.28'
$ws.Range('E51').Value = '  +0.35%  '
